$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Report Generated On" timestamp
$ws.Range("D5").Value = "Report Generated On: 08/26/2025 10:01 AM"

# Update Total Billed Amount
$ws.Range("C8").Value = 208.26

# Clear the Scope ID # value (was "#INVALID VALUE", now blank)
$ws.Range("G10").Value = ""

# Update the line item pricing and total pricing
$ws.Range("H16").Value = 208.26
$ws.Range("H17").Value = 208.26
